$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: convert the Date cell (C2) from inline-string text to a real
# date serial, formatted as YYYY-MM-DD ---
$ws.Range("C2").NumberFormat = "yyyy-mm-dd"
$ws.Range("C2").Value = "2025-07-21"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"

# --- Row 3: new progress-history entry ---
$ws.Range("A3").Value = "G1"
$ws.Range("B3").Value = "Read"

$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
$ws.Range("C3").Value = "2025-07-22"

$ws.Range("D3").Value = 1.01
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 0.01

Write-Host "done"
